$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the two reviewer comments (id 1 "This part also needs a
#    rewrite..." and id 2 "This also needs to be put into third person...").
#    The author incorporated the requested rewrites directly into the body
#    text below, so the comments (and their anchors / range markers) go
#    away entirely. Comment id 0 stays untouched.
#    Delete from the end first so earlier indices stay stable.
# ---------------------------------------------------------------------------
$d.Comments.Item(3).Delete()
$d.Comments.Item(2).Delete()

# ---------------------------------------------------------------------------
# 2) "by almost everybody." starting a brand-new sentence "Whether it is..."
#    becomes one continuous sentence: "...everybody whether it is..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "by almost everybody. Whether it is",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "by almost everybody whether it is",
    2
)

# ---------------------------------------------------------------------------
# 3) Rewrite of the "When I am out and about..." sentence (previously the
#    text wrapped by comment id 1) into the third-person-friendlier, tighter
#    version requested by the reviewer.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "When I am out and about and think about something I need to type up at a later point when I have access to my computer I can write it in to my phone and then my phone syncs to the cloud services I have and then when I turn my computer on it loads up and syncs to the cloud services and brings my saved document",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "When I am out and about and think about something I need to type into a document, I can write it on my phone and then my phone will sync to the cloud services I have, when I turn my computer on it loads up and syncs to the cloud services, bringing my saved document",
    2
)

# ---------------------------------------------------------------------------
# 4) Tidy up the closing sentence (previously the tail of the text wrapped
#    by comment id 2): add commas for readability and fix "with out" -> 
#    "without" (a single Find/Replace spanning both proofErr-wrapped runs
#    removes the now-stale spelling-error markers along with the misspelling).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "from schooling online to the games we play to the documents we share between friends family and work with out the advances",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "from schooling online, the games we play, to the documents we share between friends family and work, without the advances",
    2
)
